# Update recurrence metrics for the last row (year 2025, row 6):
# total_customers (C) and new_customers (E) each increase by 1 (one extra
# new customer), returning_customers (D) stays the same, and the derived
# rate columns new_rate (G) and returning_rate (H) are recalculated
# accordingly. retention_rate (F) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 424
$ws.Range("E6").Value = 115
$ws.Range("G6").Value = 27.12264150943397
$ws.Range("H6").Value = 72.87735849056604
